$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 3.6
$ws.Range("I2").Value = 1.91
$ws.Range("J2").Value = 3.75
$ws.Range("L2").Value = 2.5
$ws.Range("S2").Value = 1.86
$ws.Range("T2").Value = 2.04
$ws.Range("U2").Value = 2.3
$ws.Range("V2").Value = 1.62
$ws.Range("Y2").Value = 1.5
$ws.Range("Z2").Value = 2.5
$ws.Range("AE2").Value = 26
$ws.Range("AG2").Value = 19
$ws.Range("AM2").Value = 13
$ws.Range("AO2").Value = 21

# Row 4 updates
$ws.Range("G4").Value = 1.67
$ws.Range("H4").Value = 3.6
$ws.Range("I4").Value = 5.25
$ws.Range("J4").Value = 2.3
$ws.Range("K4").Value = 2.1
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.73
$ws.Range("W4").Value = 1.44
$ws.Range("X4").Value = 2.63
$ws.Range("Y4").Value = 2
$ws.Range("Z4").Value = 1.73
$ws.Range("AA4").Value = 6
$ws.Range("AC4").Value = 8.5
$ws.Range("AF4").Value = 29
$ws.Range("AG4").Value = 8.5
$ws.Range("AH4").Value = 7
$ws.Range("AK4").Value = 401
$ws.Range("AL4").Value = 12

# Row 7 updates
$ws.Range("Y7").Value = 1.8
$ws.Range("Z7").Value = 1.95
$ws.Range("AG7").Value = 15
$ws.Range("AL7").Value = 17
$ws.Range("AM7").Value = 34
